$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.810.33'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = '1.641.95'
$ws.Range('E3').Value = '  -0.31%  '
$ws.Range('E4').Value = '  -0.42%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '218.52'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('E7').Value = '  -0.39%  '
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.0623'
$ws.Range('E9').Value = '  -0.82%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.23'
$ws.Range('E10').Value = '  +0.07%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0847'
$ws.Range('E11').Value = '  +0.58%  '
$ws.Range('D12').Value = '1.870.86'
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('D13').Value = '1.666.74'
$ws.Range('E13').Value = '  +1.12%  '
$ws.Range('E14').Value = '  -0.61%  '
$ws.Range('E15').Value = '  -0.70%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.17'
$ws.Range('E16').Value = '  +0.66%  '
$ws.Range('D17').Value = '26.820.80'
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('D18').Value = '0.0₃0734'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '216.05'
$ws.Range('E19').Value = '  +0.86%  '
$ws.Range('E20').Value = '  -0.38%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.62'
$ws.Range('E21').Value = '  +4.89%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.37'
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.35'
$ws.Range('E23').Value = '  -2.25%  '
$ws.Range('E24').Value = '  -2.25%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '147.57'
$ws.Range('E25').Value = '  +1.74%  '
$ws.Range('E26').Value = '  -0.48%  '
$ws.Range('E27').Value = '  -0.17%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.13'
$ws.Range('E28').Value = '  +0.30%  '
$ws.Range('E29').Value = '  +0.13%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0509'
$ws.Range('E30').Value = '  -0.75%  '
$ws.Range('E31').Value = '  +0.80%  '
$ws.Range('E32').Value = '  +1.72%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.99'
$ws.Range('E33').Value = '  -0.79%  '
$ws.Range('E34').Value = '  +0.92%  '
$ws.Range('D35').Value = '1.266.10'
$ws.Range('E35').Value = '  -2.23%  '
$ws.Range('E36').Value = '  +0.20%  '
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('E38').Value = '  -1.72%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.818'
$ws.Range('E39').Value = '  -1.34%  '
$ws.Range('E40').Value = '  -0.28%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.803'
$ws.Range('E41').Value = '  -1.07%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.34'
$ws.Range('E42').Value = '  -0.46%  '
$ws.Range('D43').Value = '1.780.69'
$ws.Range('E43').Value = '  -0.82%  '
$ws.Range('E44').Value = '  -4.47%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '92.76'
$ws.Range('E45').Value = '  +1.17%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '61.08'
$ws.Range('E46').Value = '  +1.21%  '
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('E48').Value = '  -0.67%  '
$ws.Range('E49').Value = '  -0.67%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0965'
$ws.Range('E50').Value = '  -1.51%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.53'
$ws.Range('E51').Value = '  -2.10%  '
